$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 51231.85
$ws.Range("I6").Value = 150
$ws.Range("J6").Value = 68259.13
$ws.Range("K6").Value = 450
$ws.Range("L6").Value = 204777.39
$ws.Range("M6").Value = -338
$ws.Range("N6").Value = -205001.39
$ws.Range("H12").Value = 547.6429000000001
$ws.Range("I12").Value = 447.22223
$ws.Range("K12").Value = 447.22223
$ws.Range("M12").Value = -277.22223
$ws.Range("H58").Value = 2663.0908
$ws.Range("I58").Value = 979
$ws.Range("J58").Value = 4066.5
$ws.Range("K58").Value = 2937
$ws.Range("L58").Value = 12199.5
$ws.Range("M58").Value = -2787
$ws.Range("N58").Value = -12499.5
$ws.Range("H70").Value = 999.5
$ws.Range("H73").Value = 999.5
$ws.Range("H137").Value = 752528.2
$ws.Range("I137").Value = 982786.8
$ws.Range("K137").Value = 2948360.4
$ws.Range("M137").Value = -2945810.4
$ws.Range("H138").Value = 181877.77
$ws.Range("J138").Value = 5338.085
$ws.Range("L138").Value = 16014.255
$ws.Range("N138").Value = -26294.255
$ws.Range("H140").Value = 99999.5
$ws.Range("J140").Value = 99999.5
$ws.Range("L140").Value = 99999.5
$ws.Range("N140").Value = -110359.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 732.3333
$ws.Range("I4").Value = 348.5
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 348.5
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -232.5
$ws.Range("N4").Value = -1732
$ws.Range("H32").Value = 7478.64
$ws.Range("I32").Value = 6735.4106
$ws.Range("J32").Value = 21600
$ws.Range("K32").Value = 6735.4106
$ws.Range("L32").Value = 21600
$ws.Range("M32").Value = -6448.4106
$ws.Range("N32").Value = -22174
$ws.Range("H61").Value = 3751.3062
$ws.Range("I61").Value = 3786.4473
$ws.Range("J61").Value = 3629.9092
$ws.Range("K61").Value = 3786.4473
$ws.Range("L61").Value = 3629.9092
$ws.Range("M61").Value = -3574.4473
$ws.Range("N61").Value = -4053.9092
$ws.Range("H119").Value = 30580.4
$ws.Range("I119").Value = 21300.666
$ws.Range("J119").Value = 44500
$ws.Range("K119").Value = 21300.666
$ws.Range("L119").Value = 44500
$ws.Range("M119").Value = -16462.666
$ws.Range("N119").Value = -54176
$ws.Range("H132").Value = 2542.5386
$ws.Range("I132").Value = 1505.4667
$ws.Range("J132").Value = 5999.4443
$ws.Range("K132").Value = 4516.4001
$ws.Range("L132").Value = 17998.3329
$ws.Range("M132").Value = -1986.4001
$ws.Range("N132").Value = -23058.3329
$ws.Range("H136").Value = 3751.3062
$ws.Range("I136").Value = 3786.4473
$ws.Range("J136").Value = 3629.9092
$ws.Range("K136").Value = 11359.3419
$ws.Range("L136").Value = 10889.7276
$ws.Range("M136").Value = -8809.341899999999
$ws.Range("N136").Value = -15989.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8975.5
$ws.Range("I86").Value = 13551
$ws.Range("K86").Value = 13551
$ws.Range("M86").Value = -12428
$ws.Range("H89").Value = 8975.5
$ws.Range("I89").Value = 13551
$ws.Range("K89").Value = 67755
$ws.Range("M89").Value = -62139
$ws.Range("H99").Value = 14220.655
$ws.Range("I99").Value = 16769.408
$ws.Range("K99").Value = 16769.408
$ws.Range("M99").Value = -15271.408
$ws.Range("H105").Value = 47375.68
$ws.Range("I105").Value = 80703.78999999999
$ws.Range("K105").Value = 80703.78999999999
$ws.Range("M105").Value = -78956.78999999999
$ws.Range("H107").Value = 988.8889
$ws.Range("I107").Value = 989.625
$ws.Range("K107").Value = 989.625
$ws.Range("M107").Value = 930.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 231.6
$ws.Range("I22").Value = 199.5
$ws.Range("J22").Value = 253
$ws.Range("K22").Value = 199.5
$ws.Range("L22").Value = 253
$ws.Range("M22").Value = 150.5
$ws.Range("N22").Value = -953
$ws.Range("H31").Value = 4757.5
$ws.Range("I31").Value = 1300
$ws.Range("J31").Value = 5449
$ws.Range("K31").Value = 1300
$ws.Range("L31").Value = 5449
$ws.Range("M31").Value = -1005
$ws.Range("N31").Value = -6039
$ws.Range("H34").Value = 4757.5
$ws.Range("I34").Value = 1300
$ws.Range("J34").Value = 5449
$ws.Range("K34").Value = 1300
$ws.Range("L34").Value = 5449
$ws.Range("M34").Value = -1098
$ws.Range("N34").Value = -5853
$ws.Range("H58").Value = 1955.4857
$ws.Range("I58").Value = 2089.4546
$ws.Range("J58").Value = 1894.0834
$ws.Range("K58").Value = 2089.4546
$ws.Range("L58").Value = 1894.0834
$ws.Range("M58").Value = -1886.4546
$ws.Range("N58").Value = -2300.0834
$ws.Range("H107").Value = 5759.143
$ws.Range("I107").Value = 7743.3335
$ws.Range("J107").Value = 798.6667
$ws.Range("K107").Value = 7743.3335
$ws.Range("L107").Value = 798.6667
$ws.Range("M107").Value = -5823.3335
$ws.Range("N107").Value = -4638.6667
$ws.Range("H132").Value = 6011.4346
$ws.Range("I132").Value = 7192.0586
$ws.Range("J132").Value = 2666.3333
$ws.Range("K132").Value = 21576.1758
$ws.Range("L132").Value = 7998.999899999999
$ws.Range("M132").Value = -19046.1758
$ws.Range("N132").Value = -13058.9999
$ws.Range("H134").Value = 2066.0938
$ws.Range("I134").Value = 2189.08
$ws.Range("J134").Value = 1626.8572
$ws.Range("K134").Value = 6567.24
$ws.Range("L134").Value = 4880.571599999999
$ws.Range("M134").Value = -4032.24
$ws.Range("N134").Value = -9950.571599999999
$ws.Range("H136").Value = 1955.4857
$ws.Range("I136").Value = 2089.4546
$ws.Range("J136").Value = 1894.0834
$ws.Range("K136").Value = 6268.3638
$ws.Range("L136").Value = 5682.2502
$ws.Range("M136").Value = -3718.3638
$ws.Range("N136").Value = -10782.2502
$ws.Range("H141").Value = 409598.34
$ws.Range("J141").Value = 515148.25
$ws.Range("L141").Value = 515148.25
$ws.Range("N141").Value = -525508.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 6031.839
$ws.Range("J122").Value = 8184.7617
$ws.Range("L122").Value = 73662.8553
$ws.Range("N122").Value = -78562.8553
$ws.Range("H131").Value = 4752.6665
$ws.Range("I131").Value = 8205.727999999999
$ws.Range("J131").Value = 1830.8462
$ws.Range("K131").Value = 24617.184
$ws.Range("L131").Value = 5492.5386
$ws.Range("M131").Value = -19577.184
$ws.Range("N131").Value = -15572.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 79999
$ws.Range("J112").Value = 79999
$ws.Range("L112").Value = 79999
$ws.Range("N112").Value = -82215
$ws.Range("H132").Value = 3357.2942
$ws.Range("J132").Value = 2398.8
$ws.Range("L132").Value = 7196.400000000001
$ws.Range("N132").Value = -12256.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 29179.916
$ws.Range("I40").Value = 41494.637
$ws.Range("K40").Value = 41494.637
$ws.Range("M40").Value = -41358.637
$ws.Range("H61").Value = 14793.391
$ws.Range("I61").Value = 1841.4
$ws.Range("J61").Value = 27128.62
$ws.Range("K61").Value = 1841.4
$ws.Range("L61").Value = 27128.62
$ws.Range("M61").Value = -1639.4
$ws.Range("N61").Value = -27532.62
$ws.Range("H113").Value = 14793.391
$ws.Range("I113").Value = 1841.4
$ws.Range("J113").Value = 27128.62
$ws.Range("K113").Value = 1841.4
$ws.Range("L113").Value = 27128.62
$ws.Range("M113").Value = 328.5999999999999
$ws.Range("N113").Value = -31468.62
$ws.Range("H132").Value = 468967.1
$ws.Range("I132").Value = 995744.5600000001
$ws.Range("J132").Value = 4163.4707
$ws.Range("K132").Value = 2987233.68
$ws.Range("L132").Value = 12490.4121
$ws.Range("M132").Value = -2984703.68
$ws.Range("N132").Value = -17550.4121
$ws.Range("H136").Value = 6976.636
$ws.Range("I136").Value = 4597.4546
$ws.Range("K136").Value = 13792.3638
$ws.Range("M136").Value = -11242.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8395.273999999999
$ws.Range("I132").Value = 9981.186
$ws.Range("J132").Value = 3887.9473
$ws.Range("K132").Value = 29943.558
$ws.Range("L132").Value = 11663.8419
$ws.Range("M132").Value = -27413.558
$ws.Range("N132").Value = -16723.8419
$ws.Range("H136").Value = 15384615
$ws.Range("I136").Value = 15384615
$ws.Range("K136").Value = 46153845
$ws.Range("M136").Value = -46151295
